# Updates the Iceland MSME summary figures (revised data refresh).
# The affected cells hold their numbers as literal text (shared strings),
# so we force Text formatting before writing the new values to keep them
# stored as text rather than being auto-converted to numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B13" = "74.55"   # Enterprises density (per 1000 people) - Micro
    "C13" = "4.64"    # Enterprises density (per 1000 people) - SMEs
    "B14" = "26.31"   # Employment (% of total) - Micro
    "C14" = "45.36"   # Employment (% of total) - SMEs
    "D14" = "71.67"   # Employment (% of total) - MSMEs
    "B16" = "93.91"   # Enterprises (% of total) - Micro
    "C16" = "5.85"    # Enterprises (% of total) - SMEs
    "D16" = "99.76"   # Enterprises (% of total) - MSMEs
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
